$d = $word.ActiveDocument

# The sentence " ... sviluppatori e alcuni utenti scelti, in cui gli utenti
# provano il prodotto ..." is revised to "... in cui questi ultimi provano
# il prodotto ...". Word records the edit point with its hidden "_GoBack"
# bookmark, which moves from the end of the document (after "testing!") to
# right before the newly typed "questi ultimi".

# Step 1: split "cui " off into its own run (this must happen BEFORE the
# text substitution below so the preceding " alcuni utenti scelti, in "
# text keeps its own run instead of merging back into its neighbour).
$rCui = $d.Content
$rCui.Find.Execute("cui ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rCui.Font.Size = 13
$rCui.Font.Size = 12

# Step 2: replace "gli utenti" with "questi ultimi".
$rOld = $d.Content
$rOld.Find.Execute("gli utenti", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rOld.Text = "questi ultimi"

# Step 3: split the freshly typed "questi ultimi" into its own run so the
# following " provano il prodotto ..." text remains a separate run too.
$rNew = $d.Content
$rNew.Find.Execute("questi ultimi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rNew.Font.Size = 13
$rNew.Font.Size = 12

# Step 4: move the "_GoBack" bookmark to the collapsed point right before
# "questi ultimi" (adding a bookmark named "_GoBack" replaces/removes any
# existing one elsewhere in the document, matching Word's behaviour of
# tracking only the single most-recent edit location).
$bmRange = $d.Range($rNew.Start, $rNew.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
